# "Implementação de global index"
# Replaces the raw per-track numbers in the last couple of rows with
# shared-string "global index" aggregate figures, and appends four more
# aggregate rows (rows 40-42) below the existing data, each formatted the
# same way as the already-existing summary row (row 39 / index 37).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "global index" values per column (B..J) for rows 38..42, in the
# same left-to-right / top-to-bottom order the original sheet was built in.
$colValues = @{
    "B" = @("0.6650200000000001","0.7451333333333333","0.5243789473684211","0.7120000000000002","0.4865399999999999")
    "C" = @("0.7758200000000001","0.6195000000000002","0.8442105263157896","0.6897799999999998","0.80486")
    "D" = @("0.082386","0.07762666666666672","0.07373684210526313","0.055468","0.062732")
    "E" = @("0.06228746000000001","0.2659283333333333","0.05897290315789475","0.12148599000000004","0.0241865002")
    "F" = @("0.0470132104","0.04215750933333333","0.022514615263157896","0.13015698850000001","0.064895861")
    "G" = @("0.42242000000000013","0.6012666666666667","0.6208842105263159","0.471703","0.4537019999999999")
    "H" = @("-5.047700000000001","-7.482966666666668","-5.262652631578948","-6.581110000000003","-5.8405")
    "I" = @("121.88228","117.87040000000002","125.52905263157899","120.15134999999998","125.39779999999998")
    "J" = @("0.18999799999999997","0.14767333333333338","0.20245368421052637","0.17892300000000003","0.199264")
}

$targetRows = @(38, 39, 40, 41, 42)

# Index-column (A) values for the new rows, matching the style already used
# by A39 (bordered / centered header-ish style, s="1").
$indexValues = @{ 38 = 36; 39 = 37; 40 = 38; 41 = 39; 42 = 40 }

# Make sure rows 40-42 exist with the same look as row 39's A cell before
# writing values into them, by copying A39's formatting down.
foreach ($r in @(40, 41, 42)) {
    $ws.Range("A39").Copy() | Out-Null
    $ws.Range("A$r").PasteSpecial(-4122) | Out-Null
    $ws.Range("A$r").Value = $indexValues[$r]
}
# Row 38/39's A cell already has the value/style it needs untouched.

# Write the "global index" numbers as text-backed shared strings (full
# float precision, exactly as typed) rather than as native numeric cells:
# mark the destination as Text first so Excel doesn't re-parse the string
# back into a double, then restore the cell to the default ("Normal")
# style so no stray number format lingers on the cell.
foreach ($col in @("B","C","D","E","F","G","H","I","J")) {
    for ($i = 0; $i -lt $targetRows.Length; $i++) {
        $r = $targetRows[$i]
        $addr = "$col$r"
        $ws.Range($addr).NumberFormat = "@"
        $ws.Range($addr).Value = $colValues[$col][$i]
        $ws.Range($addr).Style = "Normal"
    }
}
